$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Insert rows for the two new FuelGroup categories ("Fossil Gases", "Fossil Liquids") ---
# Executed in ascending row order so each index is valid at the moment of insertion
# (mirrors native Excel Rows.Insert() shifting everything at/below down by one).
$ws.Rows.Item(7).Insert()
$ws.Rows.Item(10).Insert()
$ws.Rows.Item(19).Insert()
$ws.Rows.Item(22).Insert()
$ws.Rows.Item(31).Insert()
$ws.Rows.Item(34).Insert()

# --- 2) Populate the newly inserted rows ---
# Row 7: Fossil Gases (2030)
$ws.Cells.Item(7, 1).Value = "Fossil Gases"
$ws.Cells.Item(7, 2).Value = [double]"2030"
$ws.Cells.Item(7, 6).Value = [double]"0.0006702424723200279"
$ws.Cells.Item(7, 9).Value = [double]"4.254209269224784e-05"
# Row 10: Fossil Liquids (2030)
$ws.Cells.Item(10, 1).Value = "Fossil Liquids"
$ws.Cells.Item(10, 2).Value = [double]"2030"
$ws.Cells.Item(10, 6).Value = [double]"0.0310967063710098"
$ws.Cells.Item(10, 7).Value = [double]"4.467658727433266e-05"
$ws.Cells.Item(10, 8).Value = [double]"0.0006233899828972"
$ws.Cells.Item(10, 9).Value = [double]"0.0126970711122363"
$ws.Cells.Item(10, 10).Value = [double]"0.0005884652842813"
$ws.Cells.Item(10, 11).Value = [double]"0.0096127645369281"
# Row 19: Fossil Gases (2040)
$ws.Cells.Item(19, 1).Value = "Fossil Gases"
$ws.Cells.Item(19, 2).Value = [double]"2040"
$ws.Cells.Item(19, 6).Value = [double]"0.0003576098782463097"
$ws.Cells.Item(19, 9).Value = [double]"4.565330154468007e-05"
# Row 22: Fossil Liquids (2040)
$ws.Cells.Item(22, 1).Value = "Fossil Liquids"
$ws.Cells.Item(22, 2).Value = [double]"2040"
$ws.Cells.Item(22, 6).Value = [double]"0.008290632842178699"
$ws.Cells.Item(22, 7).Value = [double]"4.802981521093037e-05"
$ws.Cells.Item(22, 8).Value = [double]"0.0005877402693978"
$ws.Cells.Item(22, 9).Value = [double]"0.0057965788066297"
$ws.Cells.Item(22, 10).Value = [double]"0.0005222102927215999"
$ws.Cells.Item(22, 11).Value = [double]"0.0093257754520696"
# Row 31: Fossil Gases (2050)
$ws.Cells.Item(31, 1).Value = "Fossil Gases"
$ws.Cells.Item(31, 2).Value = [double]"2050"
$ws.Cells.Item(31, 6).Value = [double]"1.915507313125229e-05"
$ws.Cells.Item(31, 9).Value = [double]"1.624687691858186e-05"
# Row 34: Fossil Liquids (2050)
$ws.Cells.Item(34, 1).Value = "Fossil Liquids"
$ws.Cells.Item(34, 2).Value = [double]"2050"
$ws.Cells.Item(34, 6).Value = [double]"0.0004604453049551258"
$ws.Cells.Item(34, 7).Value = [double]"4.332085924059987e-05"
$ws.Cells.Item(34, 8).Value = [double]"0.0005330047601981"
$ws.Cells.Item(34, 9).Value = [double]"0.0010193311062008"
$ws.Cells.Item(34, 10).Value = [double]"0.0004494929781736"
$ws.Cells.Item(34, 11).Value = [double]"0.008694192101945299"

# --- 3) Update pre-existing rows whose figures changed now that the new categories feed them ---
# Row 9: Biogenic Liquids (2030)
$ws.Cells.Item(9, 6).Value = [double]"0.003076450763477843"
$ws.Cells.Item(9, 7).Value = [double]"6.343667119656581e-06"
$ws.Cells.Item(9, 8).Value = [double]"6.820861712839786e-05"
$ws.Cells.Item(9, 9).Value = [double]"0.001980603565556"
$ws.Cells.Item(9, 10).Value = [double]"9.721874768721148e-05"
$ws.Cells.Item(9, 11).Value = [double]"0.0009813852249140999"
# Row 13: Overall Demand (2030)
$ws.Cells.Item(13, 5).Value = [double]"0.001555883245481137"
$ws.Cells.Item(13, 6).Value = [double]"0.03506423229108502"
$ws.Cells.Item(13, 7).Value = [double]"5.102025439398923e-05"
$ws.Cells.Item(13, 8).Value = [double]"0.0006915986249547687"
$ws.Cells.Item(13, 9).Value = [double]"0.01478415980042846"
$ws.Cells.Item(13, 10).Value = [double]"0.0006856840319685115"
$ws.Cells.Item(13, 11).Value = [double]"0.0105941497618422"
# Row 21: Biogenic Liquids (2040)
$ws.Cells.Item(21, 6).Value = [double]"0.001289817512211268"
$ws.Cells.Item(21, 7).Value = [double]"1.034122395044084e-05"
$ws.Cells.Item(21, 8).Value = [double]"8.303424373222285e-05"
$ws.Cells.Item(21, 9).Value = [double]"0.0013345459759956"
$ws.Cells.Item(21, 10).Value = [double]"0.0001181049999963"
$ws.Cells.Item(21, 11).Value = [double]"0.0011103608962393"
# Row 25: Overall Demand (2040)
$ws.Cells.Item(25, 5).Value = [double]"0.002641854766313445"
$ws.Cells.Item(25, 6).Value = [double]"0.01083799728189062"
$ws.Cells.Item(25, 7).Value = [double]"5.83710391613712e-05"
$ws.Cells.Item(25, 8).Value = [double]"0.0006707765999769638"
$ws.Cells.Item(25, 9).Value = [double]"0.007274622889333377"
$ws.Cells.Item(25, 10).Value = [double]"0.0006403152927179"
$ws.Cells.Item(25, 11).Value = [double]"0.0104361363483089"
# Row 32: Synthetic Liquids (2050)
$ws.Cells.Item(32, 6).Value = [double]"7.203956199715926e-12"
$ws.Cells.Item(32, 7).Value = [double]"3.808113841955441e-13"
$ws.Cells.Item(32, 8).Value = [double]"3.199083786398123e-12"
$ws.Cells.Item(32, 9).Value = [double]"2.124473067637537e-11"
$ws.Cells.Item(32, 10).Value = [double]"8.088805598649355e-13"
$ws.Cells.Item(32, 11).Value = [double]"7.597388414361179e-11"
# Row 33: Biogenic Liquids (2050)
$ws.Cells.Item(33, 6).Value = [double]"0.0001165720561071835"
$ws.Cells.Item(33, 7).Value = [double]"1.847795157599896e-05"
$ws.Cells.Item(33, 8).Value = [double]"0.0001087569693328"
$ws.Cells.Item(33, 9).Value = [double]"0.0003461973294988"
$ws.Cells.Item(33, 10).Value = [double]"0.0001519017325245"
$ws.Cells.Item(33, 11).Value = [double]"0.0015818679863456"
# Row 37: Overall Demand (2050)
$ws.Cells.Item(37, 5).Value = [double]"0.005377804786596121"
$ws.Cells.Item(37, 6).Value = [double]"0.001769757469746319"
$ws.Cells.Item(37, 7).Value = [double]"6.179881119741021e-05"
$ws.Cells.Item(37, 8).Value = [double]"0.0006417652697674566"
$ws.Cells.Item(37, 9).Value = [double]"0.001514742912374316"
$ws.Cells.Item(37, 10).Value = [double]"0.0006013947115069806"
$ws.Cells.Item(37, 11).Value = [double]"0.01027606016426478"
